# Pin reassignment edit for fpga_pins.xlsx
# Updates the "Pin" column (F) on Sheet1 for several rows in the
# XC3S50A-VQ100 (Spartan-3A) pin table, and the Bank column (G) where needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F14").Value = "P9"
$ws.Range("F15").Value = "P71"

$ws.Range("F20").Value = "P10"

$ws.Range("F22").Value = "P70"

$ws.Range("F24").Value = "P86"
$ws.Range("G24").Value = 0

$ws.Range("F25").Value = "P5"
$ws.Range("F26").Value = "P6"

$ws.Range("F31").Value = "P83"
$ws.Range("G31").Value = 0

# Reflect the last edited / selected cell in the workbook UI state
$ws.Range("F20").Select()
